$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $TextValue) {
    $cell = $Sheet.Range($CellRef)
    # Force text storage so numeric-looking strings (e.g. "265.02")
    # are preserved verbatim as text instead of being coerced to a
    # Number cell (which would also lose fixed trailing-zero formatting).
    $cell.NumberFormat = "@"
    $cell.Value = $TextValue
    # Restore the default "Normal" style so no stray number-format
    # style is left attached to the cell.
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '265.02'
Set-TextValue $ws 'D3' '22.83'
Set-TextValue $ws 'D4' '6.226'
Set-TextValue $ws 'D5' '0.06158'
Set-TextValue $ws 'D7' '6.703'
Set-TextValue $ws 'D8' '1.361'
Set-TextValue $ws 'D9' '0.8129'
Set-TextValue $ws 'B10' 'WazirX'
Set-TextValue $ws 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D10' '0.1598'
Set-TextValue $ws 'E10' '9WazirXWRX'
Set-TextValue $ws 'B11' 'MandalaExchangeToken'
Set-TextValue $ws 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D11' '0.08202'
Set-TextValue $ws 'E11' '10MandalaExchangeTokenMDX'
Set-TextValue $ws 'B12' 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D12' '0.03372'
Set-TextValue $ws 'E12' '11LiechtensteinCryptoassetsExchangeLCX'
Set-TextValue $ws 'B13' 'BitrueCoin'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D13' '0.03140'
Set-TextValue $ws 'E13' '12BitrueCoinBTR'
Set-TextValue $ws 'B14' 'BitMartToken'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D14' '0.09239'
Set-TextValue $ws 'E14' '13BitMartTokenBMX'
Set-TextValue $ws 'B15' 'MCDex'
Set-TextValue $ws 'C15' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws 'D15' '3.899'
Set-TextValue $ws 'E15' '14MCDexMCB'
Set-TextValue $ws 'B16' 'BitForexToken'
Set-TextValue $ws 'C16' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D16' '0.001691'
Set-TextValue $ws 'E16' '15BitForexTokenBF'
Set-TextValue $ws 'B17' 'CoinExToken'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws 'D17' '0.04835'
Set-TextValue $ws 'E17' '16CoinExTokenCET'
Set-TextValue $ws 'B18' 'One'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws 'D18' '0.0006262'
Set-TextValue $ws 'E18' '17OneONEWorstin24h'
Set-TextValue $ws 'D19' '0.006192'
Set-TextValue $ws 'D20' '0.006265'
Set-TextValue $ws 'D23' '3.696'
Set-TextValue $ws 'D25' '0.3380'
Set-TextValue $ws 'D26' '0.1196'
Set-TextValue $ws 'D40' '0.04587'
Set-TextValue $ws 'D41' '0.007251'
Set-TextValue $ws 'D42' '0.1132'
Set-TextValue $ws 'D44' '0.01048'
Set-TextValue $ws 'D45' '0.00006146'
Set-TextValue $ws 'E47' '46CoinbaseStockTokenCOIN'
Set-TextValue $ws 'D48' '0.1975'

Write-Output "Applied 55 cell updates"
